$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Previously the data was written with a header row ("Key"/"Value") at row 1
# and the key/value pairs starting at row 2. The commit message states
# "added header=False on excel output. Files upload to A2" -> now the
# header row is no longer written, so the whole data block shifts up by
# one row, with the first key/value pair now starting at A1/B1.
$ws.Rows.Item(1).Delete()
